# Auto-generated edit script applying meteocat data refresh (2026-03-01 04:50 run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-03-01 04:48:29"
$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "92%"
$ws.Range("N2").Value = "-1.6 °C 4:19 TU"
$ws.Range("O2").Value = "-0.7 °C"
$ws.Range("E3").Value = "2026-03-01 04:48:31"
$ws.Range("E4").Value = "2026-03-01 04:48:34"
$ws.Range("E5").Value = "2026-03-01 04:48:37"
$ws.Range("N5").Value = "-4.4 °C 4:21 TU"
$ws.Range("O5").Value = "-3.6 °C"
$ws.Range("E6").Value = "2026-03-01 04:48:39"
$ws.Range("H6").NumberFormat = "@"
$ws.Range("H6").Value = "87%"
$ws.Range("N6").Value = "8.4 °C 4:28 TU"
$ws.Range("E7").Value = "2026-03-01 04:48:42"
$ws.Range("N7").Value = "12.9 °C 4:26 TU"
$ws.Range("E8").Value = "2026-03-01 04:48:44"
$ws.Range("N8").Value = "9.2 °C 4:24 TU"
$ws.Range("E9").Value = "2026-03-01 04:48:47"
$ws.Range("E10").Value = "2026-03-01 04:48:49"
$ws.Range("N10").Value = "5.4 °C 4:29 TU"
$ws.Range("O10").Value = "6.7 °C"
$ws.Range("E11").Value = "2026-03-01 04:48:52"
$ws.Range("N11").Value = "6.1 °C 4:17 TU"
$ws.Range("E12").Value = "2026-03-01 04:48:54"
$ws.Range("H12").NumberFormat = "@"
$ws.Range("H12").Value = "71%"
$ws.Range("E13").Value = "2026-03-01 04:48:57"
$ws.Range("H13").NumberFormat = "@"
$ws.Range("H13").Value = "91%"
$ws.Range("J13").Value = "1026.2 hPa"
$ws.Range("N13").Value = "4.1 °C 4:21 TU"
$ws.Range("E14").Value = "2026-03-01 04:48:59"
$ws.Range("H14").NumberFormat = "@"
$ws.Range("H14").Value = "98%"
$ws.Range("E15").Value = "2026-03-01 04:49:02"
$ws.Range("E16").Value = "2026-03-01 04:49:04"
$ws.Range("H16").NumberFormat = "@"
$ws.Range("H16").Value = "84%"
$ws.Range("N16").Value = "-6.0 °C 4:25 TU"
$ws.Range("O16").Value = "-4.8 °C"
$ws.Range("E17").Value = "2026-03-01 04:49:07"
$ws.Range("E18").Value = "2026-03-01 04:49:10"
$ws.Range("E19").Value = "2026-03-01 04:49:12"
$ws.Range("E20").Value = "2026-03-01 04:49:15"
$ws.Range("N20").Value = "-3.8 °C 4:26 TU"
$ws.Range("E21").Value = "2026-03-01 04:49:17"
$ws.Range("J21").Value = "1025.4 hPa"
$ws.Range("N21").Value = "6.1 °C 4:18 TU"
$ws.Range("O21").Value = "6.6 °C"
$ws.Range("E22").Value = "2026-03-01 04:49:20"
$ws.Range("L22").Value = "12.2 km/h - 295º 4:12 TU"
$ws.Range("O22").Value = "-5.2 °C"
$ws.Range("E23").Value = "2026-03-01 04:49:23"
$ws.Range("N23").Value = "-4.1 °C 4:29 TU"
$ws.Range("E24").Value = "2026-03-01 04:49:25"
$ws.Range("O24").Value = "4.3 °C"
$ws.Range("E25").Value = "2026-03-01 04:49:27"
$ws.Range("N25").Value = "-2.8 °C 4:29 TU"
$ws.Range("E26").Value = "2026-03-01 04:49:30"
$ws.Range("J26").Value = "1025.8 hPa"
$ws.Range("N26").Value = "2.4 °C 4:29 TU"
$ws.Range("E27").Value = "2026-03-01 04:49:32"
$ws.Range("N27").Value = "-1.7 °C 4:10 TU"
$ws.Range("E28").Value = "2026-03-01 04:49:35"
$ws.Range("N28").Value = "8.4 °C 4:23 TU"
$ws.Range("E29").Value = "2026-03-01 04:49:37"
$ws.Range("N29").Value = "8.7 °C 4:08 TU"
$ws.Range("O29").Value = "9.4 °C"
$ws.Range("E30").Value = "2026-03-01 04:49:40"
$ws.Range("H30").NumberFormat = "@"
$ws.Range("H30").Value = "77%"
$ws.Range("J30").Value = "1025.4 hPa"
$ws.Range("O30").Value = "10.3 °C"
$ws.Range("E31").Value = "2026-03-01 04:49:42"
$ws.Range("E32").Value = "2026-03-01 04:49:45"
$ws.Range("M32").Value = "4.2 °C 4:29 TU"
$ws.Range("O32").Value = "2.1 °C"
$ws.Range("E33").Value = "2026-03-01 04:49:47"
$ws.Range("J33").Value = "1025.6 hPa"
$ws.Range("E34").Value = "2026-03-01 04:49:49"
$ws.Range("L34").Value = "7.9 km/h - 147º 4:02 TU"
$ws.Range("N34").Value = "-0.4 °C 4:15 TU"
$ws.Range("E35").Value = "2026-03-01 04:49:52"
$ws.Range("E36").Value = "2026-03-01 04:49:54"
$ws.Range("E37").Value = "2026-03-01 04:49:57"
$ws.Range("N37").Value = "6.1 °C 4:25 TU"
$ws.Range("O37").Value = "6.3 °C"
$ws.Range("E38").Value = "2026-03-01 04:49:59"
$ws.Range("E39").Value = "2026-03-01 04:50:02"
$ws.Range("H39").NumberFormat = "@"
$ws.Range("H39").Value = "99%"
$ws.Range("L39").Value = "35.3 km/h - 203º 4:23 TU"
$ws.Range("N39").Value = "-3.6 °C 4:27 TU"
$ws.Range("E40").Value = "2026-03-01 04:50:04"
$ws.Range("N40").Value = "6.6 °C 4:29 TU"
$ws.Range("O40").Value = "7.2 °C"
$ws.Range("E41").Value = "2026-03-01 04:50:06"
$ws.Range("H41").NumberFormat = "@"
$ws.Range("H41").Value = "87%"
$ws.Range("N41").Value = "11.2 °C 4:20 TU"
$ws.Range("O41").Value = "11.7 °C"
$ws.Range("E42").Value = "2026-03-01 04:50:09"
$ws.Range("H42").NumberFormat = "@"
$ws.Range("H42").Value = "85%"
$ws.Range("N42").Value = "6.5 °C 4:25 TU"
$ws.Range("O42").Value = "8.9 °C"
$ws.Range("E43").Value = "2026-03-01 04:50:11"
$ws.Range("E44").Value = "2026-03-01 04:50:14"
$ws.Range("N44").Value = "-3.2 °C 4:21 TU"
$ws.Range("E45").Value = "2026-03-01 04:50:16"
$ws.Range("N45").Value = "3.3 °C 4:29 TU"
$ws.Range("O45").Value = "3.6 °C"
$ws.Range("E46").Value = "2026-03-01 04:50:18"
